# Fix composition plots and add correct ARF data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the observation/process error inputs (ARF data)
$ws.Range("B4").Value = 0.38
$ws.Range("B5").Value = 1.14

# Move the active selection to B6
$ws.Range("B6").Select()
